$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.382.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.54%  '
$ws.Range("D3").Value = "'2.032.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'243.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").Value = "'0.657"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'53.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.93%  '
$ws.Range("D9").Value = "'61.01"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.08%  '
$ws.Range("D10").Value = "'0.360"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.70%  '
$ws.Range("D11").Value = "'0.0738"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.20%  '
$ws.Range("E12").Value = '  -4.18%  '
$ws.Range("D13").Value = "'0.941"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.52%  '
$ws.Range("D14").Value = "'14.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.90%  '
$ws.Range("D15").Value = "'2.326.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.13%  '
$ws.Range("D16").Value = "'5.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.86%  '
$ws.Range("D17").Value = "'2.039.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.76%  '
$ws.Range("D18").Value = "'36.242.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.87%  '
$ws.Range("D19").Value = "'16.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.76%  '
$ws.Range("D20").Value = "'70.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.36%  '
$ws.Range("D21").Value = "'0.0₃0845"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.71%  '
$ws.Range("D22").Value = "'236.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").Value = "'5.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.03%  '
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("E25").Value = '  -3.66%  '
$ws.Range("E26").Value = '  -1.14%  '
$ws.Range("D27").Value = "'162.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.24%  '
$ws.Range("D28").Value = "'9.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -12.04%  '
$ws.Range("E29").Value = '  -1.57%  '
$ws.Range("E30").Value = '  -3.70%  '
$ws.Range("D31").Value = "'1.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.28%  '
$ws.Range("D32").Value = "'4.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -11.31%  '
$ws.Range("D33").Value = "'0.0590"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.94%  '
$ws.Range("D34").Value = "'4.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -10.21%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = "'0.0862"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.35%  '
$ws.Range("D37").Value = "'1.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.41%  '
$ws.Range("D38").Value = "'2.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.08%  '
$ws.Range("E39").Value = '  -7.31%  '
$ws.Range("D40").Value = "'4.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.15%  '
$ws.Range("E41").Value = '  -3.38%  '
$ws.Range("E42").Value = '  -5.15%  '
$ws.Range("D43").Value = "'1.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.84%  '
$ws.Range("D44").Value = "'92.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.34%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = "'0.0890"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.67%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = "'1.375.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.66%  '
$ws.Range("D47").Value = "'7.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +10.35%  '
$ws.Range("D48").Value = "'15.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.57%  '
$ws.Range("E49").Value = '  +1.96%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = "'2.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.59%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = "'2.217.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.98%  '
